$wb = $excel.ActiveWorkbook

# Worksheet references
$wsStartPrice = $wb.Worksheets.Item("start_price")
$wsLinear     = $wb.Worksheets.Item("Linear")
$wsNonLinear  = $wb.Worksheets.Item("NonLinear")

# --- start_price sheet: A2 value change ---
$wsStartPrice.Range("A2").Value = 6442

# --- Linear sheet: B2:B4 numeric updates, B5 autocorrelation string update ---
$wsLinear.Range("B2").Value = -0.839970414471177
$wsLinear.Range("B3").Value = -0.04858657878081703
$wsLinear.Range("B4").Value = 5248.842441826321
$wsLinear.Range("B5").Value = "[1.0, 0.1536012065535299, 0.0012439685632319988, -0.037056101619989124, -0.02865760059285038, -0.01038472944364264, 0.17359098591165462, 0.2885080160503216, 0.18793792505321497, 0.01625339526336294, -0.04764864450084378, -0.027497450824874446, -0.010946862986865324, 0.17607948099546364, 0.34925584140362564, 0.13864730807185216, -0.023296855114161033, -0.044810053048712876, -0.03918253010528676, -0.03029234865368284]"

# --- NonLinear sheet: B3:B9 numeric updates, B10 autocorrelation string update ---
$wsNonLinear.Range("B3").Value = 1.068939393939394
$wsNonLinear.Range("B4").Value = -2.605820638354445
$wsNonLinear.Range("B5").Value = -0.1295449989881216
$wsNonLinear.Range("B6").Value = 5470.890560088196
$wsNonLinear.Range("B7").Value = -3.016454347885402
$wsNonLinear.Range("B8").Value = 0.04570203215659095
$wsNonLinear.Range("B9").Value = 5011.943755968175
$wsNonLinear.Range("B10").Value = "[1.0, 0.1535762820574648, 0.0024096120938944506, -0.03441156081150587, -0.025257762548446302, -0.009542937316900492, 0.17195175349112482, 0.2850916442992038, 0.1858499185421701, 0.015623280303192136, -0.04563596799738653, -0.023817697934165532, -0.01111192415821838, 0.17654215191485642, 0.3456673874728987, 0.13571044752236783, -0.022402907311184363, -0.04331084825293218, -0.03764466228081494, -0.027885455292340638]"
